$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -887
$ws.Range("N2").Value = -1226
$ws.Range("H40").Value = 37697
$ws.Range("I40").Value = 60642.293
$ws.Range("J40").Value = 2236.0908
$ws.Range("K40").Value = 60642.293
$ws.Range("L40").Value = 2236.0908
$ws.Range("M40").Value = -60467.293
$ws.Range("N40").Value = -2586.0908
$ws.Range("H51").Value = 9579.75
$ws.Range("I51").Value = 19680.166
$ws.Range("J51").Value = 3519.5
$ws.Range("K51").Value = 19680.166
$ws.Range("L51").Value = 3519.5
$ws.Range("M51").Value = -19196.166
$ws.Range("N51").Value = -4487.5
$ws.Range("H58").Value = 490280.62
$ws.Range("I58").Value = 1165789.6
$ws.Range("J58").Value = 2413
$ws.Range("K58").Value = 3497368.8
$ws.Range("L58").Value = 7239
$ws.Range("M58").Value = -3497218.8
$ws.Range("N58").Value = -7539
$ws.Range("H74").Value = 4522.5
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 4045
$ws.Range("K74").Value = 5000
$ws.Range("L74").Value = 4045
$ws.Range("M74").Value = -4064
$ws.Range("N74").Value = -5917
$ws.Range("H77").Value = 4522.5
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 4045
$ws.Range("K77").Value = 25000
$ws.Range("L77").Value = 20225
$ws.Range("M77").Value = -20320
$ws.Range("N77").Value = -29585
$ws.Range("H80").Value = 91449.41
$ws.Range("I80").Value = 100270.9
$ws.Range("J80").Value = 84098.164
$ws.Range("K80").Value = 300812.7
$ws.Range("L80").Value = 252294.492
$ws.Range("M80").Value = -299814.7
$ws.Range("N80").Value = -254290.492
$ws.Range("H83").Value = 91449.41
$ws.Range("I83").Value = 100270.9
$ws.Range("J83").Value = 84098.164
$ws.Range("K83").Value = 902438.1
$ws.Range("L83").Value = 756883.476
$ws.Range("M83").Value = -897446.1
$ws.Range("N83").Value = -766867.476
$ws.Range("H87").Value = 31319
$ws.Range("I87").Value = 30800
$ws.Range("J87").Value = 31838
$ws.Range("K87").Value = 30800
$ws.Range("L87").Value = 31838
$ws.Range("M87").Value = -29552
$ws.Range("N87").Value = -34334
$ws.Range("H90").Value = 31319
$ws.Range("I90").Value = 30800
$ws.Range("J90").Value = 31838
$ws.Range("K90").Value = 92400
$ws.Range("L90").Value = 95514
$ws.Range("M90").Value = -86160
$ws.Range("N90").Value = -107994
$ws.Range("H138").Value = 2404.8235
$ws.Range("I138").Value = 2359.3572
$ws.Range("J138").Value = 2413.7888
$ws.Range("K138").Value = 7078.071599999999
$ws.Range("L138").Value = 7241.366399999999
$ws.Range("M138").Value = -1938.071599999999
$ws.Range("N138").Value = -17521.3664

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30168.166
$ws.Range("I32").Value = 4919.939
$ws.Range("J32").Value = 142637.55
$ws.Range("K32").Value = 4919.939
$ws.Range("L32").Value = 142637.55
$ws.Range("M32").Value = -4632.939
$ws.Range("N32").Value = -143211.55

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2652.4583
$ws.Range("I134").Value = 2716.55
$ws.Range("J134").Value = 2332
$ws.Range("K134").Value = 8149.650000000001
$ws.Range("L134").Value = 6996
$ws.Range("M134").Value = -5614.650000000001
$ws.Range("N134").Value = -12066

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H104").Value = 32530
$ws.Range("J104").Value = 32530
$ws.Range("L104").Value = 32530
$ws.Range("N104").Value = -37772
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("H132").Value = 46879076
$ws.Range("I132").Value = 47623276
$ws.Range("J132").Value = 45458336
$ws.Range("K132").Value = 142869828
$ws.Range("L132").Value = 136375008
$ws.Range("M132").Value = -142867298
$ws.Range("N132").Value = -136380068
$ws.Range("N108").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 9812.521000000001
$ws.Range("J5").Value = 14415.6
$ws.Range("L5").Value = 43246.8
$ws.Range("N5").Value = -43470.8
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("H32").Value = 5900
$ws.Range("J32").Value = 9800
$ws.Range("L32").Value = 29400
$ws.Range("N32").Value = -29966
$ws.Range("H74").Value = 3499.5
$ws.Range("I74").Value = 1999
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 5997
$ws.Range("L74").Value = 15000
$ws.Range("M74").Value = -4936
$ws.Range("N74").Value = -17122
$ws.Range("H77").Value = 3499.5
$ws.Range("I77").Value = 1999
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 17991
$ws.Range("L77").Value = 45000
$ws.Range("M77").Value = -12687
$ws.Range("N77").Value = -55608
$ws.Range("H113").Value = 817.7143
$ws.Range("I113").Value = 850
$ws.Range("J113").Value = 804.8
$ws.Range("K113").Value = 2550
$ws.Range("L113").Value = 2414.4
$ws.Range("M113").Value = -380
$ws.Range("N113").Value = -6754.4
$ws.Range("H131").Value = 716.78
$ws.Range("J131").Value = 777.1585
$ws.Range("L131").Value = 2331.4755
$ws.Range("N131").Value = -12411.4755
$ws.Range("H135").Value = 9812.521000000001
$ws.Range("J135").Value = 14415.6
$ws.Range("L135").Value = 129740.4
$ws.Range("N135").Value = -134810.4
$ws.Range("N16").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1682651
$ws.Range("I126").Value = 2157.6428
$ws.Range("J126").Value = 2802979.8
$ws.Range("K126").Value = 6472.928400000001
$ws.Range("L126").Value = 8408939.399999999
$ws.Range("M126").Value = -4002.928400000001
$ws.Range("N126").Value = -8413879.399999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2275.3333
$ws.Range("I7").Value = 1615.5385
$ws.Range("K7").Value = 1615.5385
$ws.Range("M7").Value = -1503.5385
$ws.Range("H22").Value = 1000.1667
$ws.Range("I22").Value = 1150
$ws.Range("J22").Value = 970.2
$ws.Range("K22").Value = 1150
$ws.Range("L22").Value = 970.2
$ws.Range("M22").Value = -855
$ws.Range("N22").Value = -1560.2
$ws.Range("H27").Value = 1000.1667
$ws.Range("I27").Value = 1150
$ws.Range("J27").Value = 970.2
$ws.Range("K27").Value = 1150
$ws.Range("L27").Value = 970.2
$ws.Range("M27").Value = -1043
$ws.Range("N27").Value = -1184.2
$ws.Range("H40").Value = 57183.777
$ws.Range("I40").Value = 126300.5
$ws.Range("J40").Value = 1890.4
$ws.Range("K40").Value = 126300.5
$ws.Range("L40").Value = 1890.4
$ws.Range("M40").Value = -126164.5
$ws.Range("N40").Value = -2162.4
$ws.Range("H126").Value = 2275.3333
$ws.Range("I126").Value = 1615.5385
$ws.Range("K126").Value = 4846.6155
$ws.Range("M126").Value = -2376.6155

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 286512.84
$ws.Range("I81").Value = 250422.75
$ws.Range("J81").Value = 334633
$ws.Range("K81").Value = 500845.5
$ws.Range("L81").Value = 669266
$ws.Range("M81").Value = -499784.5
$ws.Range("N81").Value = -671388
$ws.Range("H84").Value = 286512.84
$ws.Range("I84").Value = 250422.75
$ws.Range("J84").Value = 334633
$ws.Range("K84").Value = 2504227.5
$ws.Range("L84").Value = 3346330
$ws.Range("M84").Value = -2498923.5
$ws.Range("N84").Value = -3356938
$ws.Range("H126").Value = 2490
$ws.Range("I126").Value = 3980
$ws.Range("K126").Value = 11940
$ws.Range("M126").Value = -9470
